$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header for the "order_position" column (G). Plain text, so a normal
# value assignment is enough - Excel will not try to reinterpret it as a
# number.
$ws.Range("G1").Value = "order_position"

# The order_position values themselves look numeric ("1.10", "2.5", ...),
# so a direct .Value assignment would be auto-coerced by Excel into a
# number (dropping the significant trailing zero on "1.10", etc). To keep
# them as genuine text - matching the source data - write them as literal
# text formulas first, then convert the formulas to static values via
# copy / paste-special (values only). This preserves the string content
# exactly without touching any cell's NumberFormat/Style.
$ws.Range("G2").Formula = '="1.10"'
$ws.Range("G3").Formula = '="2.5"'
$ws.Range("G4").Formula = '="1.4"'
$ws.Range("G5").Formula = '="2.6"'
$ws.Range("G6").Formula = '="3.7"'
$ws.Range("G7").Formula = '="8.9"'
$ws.Range("G8").Formula = '="10.11"'
$ws.Range("G9").Formula = '="12.13"'
$ws.Range("G10").Formula = '="14.15"'
$ws.Range("G11").Formula = '="3.7"'
$ws.Range("G12").Formula = '="4.8"'

$numRng = $ws.Range("G2:G12")
$numRng.Copy()
$numRng.PasteSpecial(-4163)

# These remaining rows use the existing "none" label, which stays text on
# a plain assignment since it doesn't parse as a number.
$ws.Range("G13").Value = "none"
$ws.Range("G14").Value = "none"
$ws.Range("G15").Value = "none"
$ws.Range("G16").Value = "none"
$ws.Range("G17").Value = "none"

# Column G width, matching the bestFit/customWidth seen in the saved file
# (the runtime quantizes ColumnWidth to 1/6-character steps, so 13.5 is the
# input that lands closest to the source file's 14.28515625).
$ws.Columns.Item(7).ColumnWidth = 13.5

# Move the active selection, matching the saved workbook view.
$ws.Range("F21").Select()
